# Update "想去人数" (want-to-go count) figures in the "展览" and "全部类型" sheets
# to reflect newly generated output (gh-pages rebuild).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (exhibitions) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 109
$ws1.Range("F3").Value = 12117
$ws1.Range("F4").Value = 49
$ws1.Range("F5").Value = 232
$ws1.Range("F6").Value = 371
$ws1.Range("F8").Value = 12025
$ws1.Range("F10").Value = 1185
$ws1.Range("F12").Value = 593
$ws1.Range("F13").Value = 1802
$ws1.Range("F14").Value = 5937

# --- Sheet "全部类型" (all types) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 109
$ws4.Range("F5").Value = 12117
$ws4.Range("F6").Value = 49
$ws4.Range("F7").Value = 232
$ws4.Range("F9").Value = 371
$ws4.Range("F11").Value = 12025
$ws4.Range("F13").Value = 1185
$ws4.Range("F15").Value = 593
$ws4.Range("F16").Value = 1802
$ws4.Range("F18").Value = 5937
